$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new 2022 column (S) data, mirroring column R (2021)
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Copy formatting from column R to column S
$ws.Range("R2:R6").Copy()
$ws.Range("S2:S6").PasteSpecial(-4122) # xlPasteFormats

# Update the selected cell in the sheet view
$ws.Range("C19").Select()
